# Apply crypto price/volume updates per commit diff (Wed Aug 23 23:43:23 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'26.496.56"
$ws.Cells.Item(2, 5).Value = "'  +1.88%  "
$ws.Cells.Item(3, 4).Value = "'1.681.17"
$ws.Cells.Item(3, 5).Value = "'  +2.89%  "
$ws.Cells.Item(4, 5).Value = "'  +0.12%  "
$ws.Cells.Item(5, 4).Value = "'216.81"
$ws.Cells.Item(5, 5).Value = "'  +2.82%  "
$ws.Cells.Item(6, 5).Value = "'  +1.72%  "
$ws.Cells.Item(7, 4).Value = "'1.002"
$ws.Cells.Item(7, 5).Value = "'  +0.14%  "
$ws.Cells.Item(8, 4).Value = "'0.2686"
$ws.Cells.Item(8, 5).Value = "'  +3.94%  "
$ws.Cells.Item(9, 4).Value = "'0.06402"
$ws.Cells.Item(9, 5).Value = "'  +1.95%  "
$ws.Cells.Item(10, 4).Value = "'21.69"
$ws.Cells.Item(10, 5).Value = "'  +5.38%  "
$ws.Cells.Item(11, 4).Value = "'0.07797"
$ws.Cells.Item(11, 5).Value = "'  +2.75%  "
$ws.Cells.Item(12, 4).Value = "'1.686.38"
$ws.Cells.Item(12, 5).Value = "'  +3.41%  "
$ws.Cells.Item(13, 4).Value = "'4.498"
$ws.Cells.Item(13, 5).Value = "'  +1.79%  "
$ws.Cells.Item(14, 4).Value = "'0.5573"
$ws.Cells.Item(14, 5).Value = "'  +1.38%  "
$ws.Cells.Item(15, 4).Value = "'0.0₅8328"
$ws.Cells.Item(15, 5).Value = "'  +3.81%  "
$ws.Cells.Item(16, 4).Value = "'65.63"
$ws.Cells.Item(16, 5).Value = "'  +1.28%  "
$ws.Cells.Item(17, 4).Value = "'26.543.04"
$ws.Cells.Item(17, 5).Value = "'  +2.14%  "
$ws.Cells.Item(18, 5).Value = "'  +0.02%  "
$ws.Cells.Item(19, 4).Value = "'4.758"
$ws.Cells.Item(19, 5).Value = "'  +1.74%  "
$ws.Cells.Item(20, 4).Value = "'194.50"
$ws.Cells.Item(20, 5).Value = "'  +4.86%  "
$ws.Cells.Item(21, 5).Value = "'  +2.15%  "
$ws.Cells.Item(22, 4).Value = "'6.353"
$ws.Cells.Item(23, 5).Value = "'  +0.12%  "
$ws.Cells.Item(24, 4).Value = "'143.10"
$ws.Cells.Item(24, 5).Value = "'  -1.52%  "
$ws.Cells.Item(25, 4).Value = "'0.1282"
$ws.Cells.Item(25, 5).Value = "'  +5.73%  "
$ws.Cells.Item(26, 4).Value = "'7.436"
$ws.Cells.Item(26, 5).Value = "'  +0.66%  "
$ws.Cells.Item(27, 4).Value = "'16.33"
$ws.Cells.Item(27, 5).Value = "'  +4.30%  "
$ws.Cells.Item(28, 5).Value = "'  +4.02%  "
$ws.Cells.Item(29, 4).Value = "'0.06232"
$ws.Cells.Item(29, 5).Value = "'  +5.40%  "
$ws.Cells.Item(30, 5).Value = "'  +2.52%  "
$ws.Cells.Item(31, 5).Value = "'  +5.22%  "
$ws.Cells.Item(32, 5).Value = "'  +1.77%  "
$ws.Cells.Item(33, 4).Value = "'1.692"
$ws.Cells.Item(33, 5).Value = "'  +4.19%  "
$ws.Cells.Item(34, 4).Value = "'1.009"
$ws.Cells.Item(34, 5).Value = "'  +2.85%  "
$ws.Cells.Item(35, 4).Value = "'2.426"
$ws.Cells.Item(35, 5).Value = "'  +1.72%  "
$ws.Cells.Item(36, 5).Value = "'  +1.50%  "
$ws.Cells.Item(37, 4).Value = "'0.5736"
$ws.Cells.Item(37, 5).Value = "'  -0.86%  "
$ws.Cells.Item(38, 4).Value = "'0.01639"
$ws.Cells.Item(39, 4).Value = "'6.034"
$ws.Cells.Item(39, 5).Value = "'  +6.41%  "
$ws.Cells.Item(40, 4).Value = "'1.074.46"
$ws.Cells.Item(40, 5).Value = "'  +3.51%  "
$ws.Cells.Item(41, 4).Value = "'0.8596"
$ws.Cells.Item(41, 5).Value = "'  +1.22%  "
$ws.Cells.Item(42, 5).Value = "'  -0.35%  "
$ws.Cells.Item(43, 4).Value = "'100.03"
$ws.Cells.Item(43, 5).Value = "'  -0.11%  "
$ws.Cells.Item(44, 4).Value = "'1.826.89"
$ws.Cells.Item(44, 5).Value = "'  +2.53%  "
$ws.Cells.Item(45, 2).Value = "'Aave"
$ws.Cells.Item(45, 3).Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(45, 4).Value = "'57.12"
$ws.Cells.Item(45, 5).Value = "'  +3.99%  "
$ws.Cells.Item(46, 2).Value = "'BabyDogeCoin"
$ws.Cells.Item(46, 3).Value = "'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Cells.Item(46, 4).Value = "'0.0₈108"
$ws.Cells.Item(46, 5).Value = "'  -1.07%  "
$ws.Cells.Item(47, 2).Value = "'Frax"
$ws.Cells.Item(47, 3).Value = "'https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Cells.Item(47, 4).Value = "'1.002"
$ws.Cells.Item(47, 5).Value = "'  +0.61%  "
$ws.Cells.Item(48, 2).Value = "'EnergySwap"
$ws.Cells.Item(48, 3).Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(48, 4).Value = "'8.115"
$ws.Cells.Item(48, 5).Value = "'  +1.10%  "
$ws.Cells.Item(49, 4).Value = "'0.05210"
$ws.Cells.Item(49, 5).Value = "'  +0.91%  "
$ws.Cells.Item(50, 2).Value = "'Aptos"
$ws.Cells.Item(50, 3).Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(50, 4).Value = "'6.029"
$ws.Cells.Item(50, 5).Value = "'  +3.05%  "
$ws.Cells.Item(51, 2).Value = "'Mantle"
$ws.Cells.Item(51, 3).Value = "'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Cells.Item(51, 4).Value = "'0.4241"
$ws.Cells.Item(51, 5).Value = "'  +0.43%  "
